# Swap the observation data between paired rows (5<->6, 7<->8, 13<->14, 18<->19).
# Each pair has identical "static" columns (locality, county, dates, observer, ...)
# but the diff shows the per-observation columns (A, B, D, E, F, G, H, Q, R, AC)
# being exchanged between the two rows of each pair.

function Swap-RowData {
    param($ws, $row1, $row2)

    $cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "AC")

    foreach ($col in $cols) {
        $addr1 = "$col$row1"
        $addr2 = "$col$row2"
        $val1 = $ws.Range($addr1).Value2
        $val2 = $ws.Range($addr2).Value2

        if ($null -eq $val2) {
            $ws.Range($addr1).Value = ""
        } else {
            $ws.Range($addr1).Value = $val2
        }

        if ($null -eq $val1) {
            $ws.Range($addr2).Value = ""
        } else {
            $ws.Range($addr2).Value = $val1
        }
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Swap-RowData $ws 5 6
Swap-RowData $ws 7 8
Swap-RowData $ws 13 14
Swap-RowData $ws 18 19
